$d = $word.ActiveDocument

# --------------------------------------------------------------------
# 1. Rename the "References" heading to "Bibliographie" and rename its
#    bookmark from "references" to "bibliographie" accordingly.
#    (Done first, while the Paragraphs collection is still reliable --
#    mutating the table's Title later can desync paragraph lookups.)
# --------------------------------------------------------------------
$headingPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "References`r") {
        $headingPara = $p
        break
    }
}

$newHeadingXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr><w:pStyle w:val="Heading1"/></w:pPr>' +
    '<w:bookmarkStart w:id="25" w:name="bibliographie"/>' +
    '<w:r><w:t xml:space="preserve">Bibliographie</w:t></w:r>' +
    '<w:bookmarkEnd w:id="25"/>' +
    '</w:p>'

$headingPara.Range.InsertXML($newHeadingXml)

# --------------------------------------------------------------------
# 2. Fix the table-caption paragraph text:
#    "Table 1: Tabeau de paramètres ..." -> "Table 1: Paramètres ..."
# --------------------------------------------------------------------
$d.Content.Find.Execute(
    "Table 1: Tabeau de paramètres utilisés pour construitre le modèle 1.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Table 1: Paramètres utilisés pour construitre le modèle 1.", 2)

# --------------------------------------------------------------------
# 3. Fix the table's accessibility title (serialised as <w:tblCaption>),
#    which holds the same caption text (with a trailing space).
# --------------------------------------------------------------------
$t = $d.Tables.Item(1)
$t.Title = "Table 1: Paramètres utilisés pour construitre le modèle 1. "
